# Inclusão de coluna numérica, representando notas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column C
$ws.Range("C1").Value = "Note_num"

# Numeric values for C2:C100 (Note_num), row-aligned with existing data in A/B
$noteNums = @(0,2,5,7,10,11,0,3,4,3,9,9,9,9,7,5,2,0,5,11,11,5,5,7,7,9,4,5,4,2,0,2,4,5,7,9,11,11,0,2,4,5,4,3,0,9,0,11,4,2,7,1,5,4,0,0,0,0,9,4,5,9,11,0,0,1,11,10,5,6,11,1,5,4,4,5,4,11,11,11,10,8,10,0,8,3,3,2,0,3,9,1,4,5,5,11,9,7,8)

for ($i = 0; $i -lt $noteNums.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $noteNums[$i]
}

# Mirror the final selection seen in the saved workbook
$ws.Range("C101").Select()
